$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, preserving the default (unstyled)
# cell format. Forcing NumberFormat "@" before the assignment stops Excel
# from reinterpreting numeric-looking strings (e.g. "612.00", "43.80") as
# numbers and dropping trailing zeros / decimal formatting; resetting the
# style back to "Normal" afterwards removes the now-unneeded text format
# so the cell keeps style index 0, same as in the source workbook.
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "92.338.67"
Set-TextValue "E2" "  +0.85%  "
Set-TextValue "D3" "3.106.16"
Set-TextValue "E3" "  -0.71%  "
Set-TextValue "E4" "  +0.01%  "
Set-TextValue "D5" "234.11"
Set-TextValue "E5" "  -3.05%  "
Set-TextValue "D6" "612.00"
Set-TextValue "E6" "  -1.05%  "
Set-TextValue "E7" "  -2.77%  "
Set-TextValue "E8" "  -0.40%  "
Set-TextValue "E9" "  -0.02%  "
Set-TextValue "D10" "3.103.41"
Set-TextValue "E10" "  -0.82%  "
Set-TextValue "D11" "0.781"
Set-TextValue "E11" "  +4.05%  "
Set-TextValue "E12" "  -3.23%  "
Set-TextValue "E13" "  -4.46%  "
Set-TextValue "D14" "92.122.18"
Set-TextValue "E14" "  +0.78%  "
Set-TextValue "B15" "Avalanche"
Set-TextValue "C15" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D15" "33.76"
Set-TextValue "E15" "  -4.07%  "
Set-TextValue "B16" "Toncoin"
Set-TextValue "C16" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D16" "5.41"
Set-TextValue "E16" "  -3.51%  "
Set-TextValue "D17" "3.687.52"
Set-TextValue "E17" "  -0.77%  "
Set-TextValue "D18" "3.070.35"
Set-TextValue "E18" "  -2.89%  "
Set-TextValue "E19" "  -0.04%  "
Set-TextValue "D20" "14.56"
Set-TextValue "E20" "  -2.84%  "
Set-TextValue "D21" "5.78"
Set-TextValue "E21" "  -2.24%  "
Set-TextValue "E22" "  +0.99%  "
Set-TextValue "D23" "9.23"
Set-TextValue "E23" "  +0.08%  "
Set-TextValue "D24" "437.54"
Set-TextValue "E24" "  -4.35%  "
Set-TextValue "D26" "85.24"
Set-TextValue "E26" "  -4.20%  "
Set-TextValue "D27" "11.45"
Set-TextValue "E27" "  -2.95%  "
Set-TextValue "D28" "3.267.83"
Set-TextValue "E28" "  -1.34%  "
Set-TextValue "E29" "  -0.03%  "
Set-TextValue "E30" "  +5.16%  "
Set-TextValue "D31" "0.229"
Set-TextValue "E31" "  -1.24%  "
Set-TextValue "E32" "  +45.29%  "
Set-TextValue "E33" "  -22.55%  "
Set-TextValue "D34" "9.16"
Set-TextValue "E34" "  -2.07%  "
Set-TextValue "D35" "7.98"
Set-TextValue "E35" "  +6.95%  "
Set-TextValue "D36" "0.156"
Set-TextValue "E36" "  -10.46%  "
Set-TextValue "D37" "25.74"
Set-TextValue "E37" "  -2.61%  "
Set-TextValue "D38" "3.97"
Set-TextValue "E38" "  +0.69%  "
Set-TextValue "E39" "  -3.28%  "
Set-TextValue "E40" "  +7.71%  "
Set-TextValue "D41" "1.28"
Set-TextValue "E41" "  -3.10%  "
Set-TextValue "D42" "464.85"
Set-TextValue "E42" "  -5.61%  "
Set-TextValue "E43" "  -1.95%  "
Set-TextValue "E44" "  -4.05%  "
Set-TextValue "D46" "160.04"
Set-TextValue "E46" "  +2.33%  "
Set-TextValue "D47" "0.682"
Set-TextValue "E47" "  -3.66%  "
Set-TextValue "E48" "  -5.07%  "
Set-TextValue "E49" "  -0.42%  "
Set-TextValue "E50" "  -3.20%  "
Set-TextValue "D51" "43.80"
Set-TextValue "E51" "  -0.48%  "
